# fixed signs on water quality parameters
#
# The sheet's header row had its influence-source columns re-ordered, and
# every data row underneath had its "rural communities to gw quality"
# value (column B) and "small growers to gw quality" value (column F)
# swapped with each other - correcting which column each figure's +/-1
# sign had ended up in.
#
# Header row 1: columns B..F rotate one step to the right (old F's header
# moves into B; old B,C,D,E each shift one column right into C,D,E,F).
#
# Data rows 2-28: column B and column F simply swap places (value AND
# cell style/format travel together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch column, far outside the used range (A1:H28), used as a holding
# slot while two cells trade places. Copy($dest) is the only form that
# carries the cell's style along with its value in this host, so a swap
# is done via a temp holder rather than plain value assignment.
$scratchCol = 26

function Swap-Cells($row, $col1, $col2) {
    $c1 = $ws.Cells.Item($row, $col1)
    $c2 = $ws.Cells.Item($row, $col2)
    $tmp = $ws.Cells.Item($row, $scratchCol)

    $c1.Copy($tmp)
    $c1.ClearContents()
    $c2.Copy($c1)
    $c2.ClearContents()
    $tmp.Copy($c2)
    $tmp.Clear()
}

# --- Header row 1: rotate B,C,D,E,F one step to the right (F wraps to B) ---
$b1 = $ws.Cells.Item(1, 2)
$c1h = $ws.Cells.Item(1, 3)
$d1h = $ws.Cells.Item(1, 4)
$e1h = $ws.Cells.Item(1, 5)
$f1 = $ws.Cells.Item(1, 6)
$tmp1 = $ws.Cells.Item(1, $scratchCol)

$f1.Copy($tmp1)   # tmp1 = old F1
$e1h.Copy($f1)    # F1   = old E1
$d1h.Copy($e1h)   # E1   = old D1
$c1h.Copy($d1h)   # D1   = old C1
$b1.Copy($c1h)    # C1   = old B1
$tmp1.Copy($b1)   # B1   = old F1
$tmp1.Clear()

# --- Data rows 2-28: swap column B and column F (value + style) ---
for ($r = 2; $r -le 28; $r++) {
    Swap-Cells $r 2 6
}

# --- Restore the active selection in the bottom-right (frozen) pane ---
$ws.Range("I3").Select()
